$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in RECORD_DUPLICATE description (row 5, column C):
# "potvrdením" -> "potvrdeným"
$ws.Range("C5").Value = "Záznam je potvrdeným duplikátom už existujúceho záznamu"

# Insert two new rows before the final "200 / SUCESS" row (currently row 14),
# pushing it down to row 16.
$ws.Range("A14:A15").EntireRow.Insert()

# New row 14: 112 / DECODING_FAILURE
$ws.Range("A14").Value = 112
$ws.Range("B14").Value = "DECODING_FAILURE"
$ws.Range("C14").Value = "Pre dopytovaný súbor nebolo možné nájsť enkódovanie. Testované sú formáty utf-8"

# New row 15: 113 / UNSUPPORTED_LOG
$ws.Range("A15").Value = 113
$ws.Range("B15").Value = "UNSUPPORTED_LOG"
$ws.Range("C15").Value = "Zvolený súbor nie je podporovaný log (PAP, KAM resp. PAP a KAM)"
$ws.Range("D15").Value = "Skontrolujte, že názov súboru obsahuje `"KAM`" alebo `"PAP`" a príponu .log"

# Grow the table (Table1) to cover the two new rows.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:D16"))

# Column D now holds longer text (the new UNSUPPORTED_LOG / "KAM"/"PAP" row),
# so re-fit its width to the new content, like Excel does automatically.
$ws.Columns.Item(4).ColumnWidth = 66.7

# Update selection to match the target workbook state.
$ws.Range("D10").Select()
